$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Ty Jerome"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Cleveland Cavaliers"

$ws.Range("A9").Value = "Anthony Davis"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Los Angeles Lakers"

$ws.Range("A11").Value = "Isaiah Hartenstein"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Oklahoma City Thunder"

$ws.Range("A12").Value = "Cameron Johnson"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Brooklyn Nets"

$ws.Range("A13").Value = "Naz Reid"
$ws.Range("B13").Value = "PF,C"
$ws.Range("C13").Value = "Minnesota Timberwolves"

$ws.Range("A16").Value = "Cade Cunningham"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Detroit Pistons"

$ws.Range("A17").Value = "Cam Thomas"
$ws.Range("B17").Value = "SG,SF"
$ws.Range("C17").Value = "Brooklyn Nets"

$ws.Range("A18").Value = "Brandon Ingram"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "New Orleans Pelicans"
